$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.259.83"
$ws.Range("E2").Value = "  +2.03%  "

$ws.Range("D3").Value = "1.874.59"
$ws.Range("E3").Value = "  +4.47%  "

$ws.Range("E4").Value = "  +0.26%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.20"
$ws.Range("E5").Value = "  +1.94%  "

$ws.Range("E6").Value = "  +0.24%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5063"
$ws.Range("E7").Value = "  +2.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3942"
$ws.Range("E8").Value = "  +2.50%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09584"
$ws.Range("E9").Value = "  +2.34%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.145"
$ws.Range("E10").Value = "  +5.10%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.83"
$ws.Range("E11").Value = "  +0.93%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.488"
$ws.Range("E12").Value = "  +3.79%  "

$ws.Range("E13").Value = "  +3.27%  "

$ws.Range("D14").Value = "1.880.69"
$ws.Range("E14").Value = "  +4.70%  "

$ws.Range("B15").Value = "BinanceUSD"
$ws.Range("C15").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.004"
$ws.Range("E15").Value = "  +0.36%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.427"
$ws.Range("E16").Value = "  +3.94%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001132"
$ws.Range("E17").Value = "  +2.43%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.96"
$ws.Range("E18").Value = "  +1.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06602"
$ws.Range("E19").Value = "  +0.99%  "

$ws.Range("E20").Value = "  +3.70%  "

$ws.Range("E21").Value = "  +0.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.186"
$ws.Range("E22").Value = "  +5.15%  "

$ws.Range("D23").Value = "28.320.52"
$ws.Range("E23").Value = "  +2.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.30"
$ws.Range("E24").Value = "  +3.66%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.304"
$ws.Range("E25").Value = "  +3.57%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.569"
$ws.Range("E26").Value = "  +7.32%  "

$ws.Range("D27").Value = "2.096.90"
$ws.Range("E27").Value = "  +4.73%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.21"
$ws.Range("E28").Value = "  +4.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "158.86"
$ws.Range("E29").Value = "  +1.47%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.61"
$ws.Range("E30").Value = "  +1.23%  "

$ws.Range("E31").Value = "  +0.31%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.067"
$ws.Range("E32").Value = "  +1.66%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.647"
$ws.Range("E33").Value = "  +2.68%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.626"
$ws.Range("E34").Value = "  +0.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.558"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06717"
$ws.Range("E36").Value = "  -1.04%  "

$ws.Range("E37").Value = "  +4.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2193"
$ws.Range("E38").Value = "  +3.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6380"
$ws.Range("E39").Value = "  +4.59%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.51"
$ws.Range("E40").Value = "  +1.75%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.004"
$ws.Range("E41").Value = "  +2.38%  "

$ws.Range("E42").Value = "  +4.72%  "

$ws.Range("E43").Value = "  +0.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.55"
$ws.Range("E44").Value = "  +4.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5992"
$ws.Range("E45").Value = "  +2.71%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.661"
$ws.Range("E46").Value = "  +0.04%  "

$ws.Range("E47").Value = "  -0.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.004"
$ws.Range("E48").Value = "  +4.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.18"
$ws.Range("E49").Value = "  +1.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06859"
$ws.Range("E51").Value = "  +2.49%  "
